$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new player row (Khris Middleton, SF, Milwaukee Bucks) was added at the
# bottom of the roster table, immediately after row 18 (Kawhi Leonard).
$ws.Range("A19").Value = "Khris Middleton"
$ws.Range("B19").Value = "SF"
$ws.Range("C19").Value = "Milwaukee Bucks"
